$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The existing last paragraph ("27/04/2014 1hr. ... Unknown".)
#    gains proofErr spellStart/spellEnd markers around the second
#    "Unknown" run, and loses its trailing _GoBack bookmark (which
#    moves to the newly-added paragraph below).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

$para1Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">27/04/2014 1hr. Crea CCD a partir de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fwrite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> usando el nombre del file como nombre y el primer parámetro como nombre de atributo, dejando &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Unknown</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221; el tipo del atributo. Por el momento, Las CCD a partir de escritura de archivo, son creadas únicamente si las llamadas a dichas funciones son realizadas afuera de cualquier flujo de control, pero ambas están preparadas para guardar atributos utilizados adentro de flujo de control, dejando como tipo del atributo &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Unknown</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;.</w:t></w:r></w:p>
"@

$null = $lastRange.InsertXML($para1Xml)

# ------------------------------------------------------------------
# 2) Append a blank paragraph plus the new "29/04/2014" log entry
#    paragraph (which now owns the _GoBack bookmark).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$null = $lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$para2Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">29/04/2014 1hr. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Parsing approach changed. Now, it allows parsing, not only from a file but also from any Iterator&lt;Character&gt;,in this way, parser scope is broadened and its logic can be used also to parse blocks inside the original parsing.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Me permite utilizar</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tokenListFactory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para los control </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>flow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> blocks adentro del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para las CCD de escritura de archivos.</w:t></w:r></w:p>
"@

$null = $newRange.InsertXML($para2Xml)

Write-Output "done"
